$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("D2").Formula = "=SUM(INDIRECT(""SINUS""))"
$ws1.Range("D3").Formula = "=SQRT(D2)"
$ws1.Range("D4").Formula = "=INDIRECT(""A1:A1"")"

$ws2.Range("A1").Formula = "=SUM(INDIRECT(""Sheet1!A1:B18""))"

$ws1.Range("A1").Select() | Out-Null

